# LOM3232.xlsx edit: 2022-09-26 16:07:08 UTC build
#
# The "Ficha de disciplina" sheet had several of its text fields swapped out
# (Objetivos / Programa resumido / Programa / Método / Bibliografia bodies
# replaced or removed) and the trailing "Requisitos" value row was removed,
# shrinking the used range from A1:C24 to A1:C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlPasteValues = -4163

# ---- Row 10 ("Objetivos:") - long body text replaced with the professor credential line.
$ws.Range("B10").Value = "5840793 - Sérgio Schneider"
$ws.Range("C10").Value = "5840793 - Sérgio Schneider"

# ---- Row 13 used to hold only "5840793 - Sérgio Schneider" (no label in A).
# It now carries the "Programa resumido:" label plus a short "Semestral" value.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A14").Copy()
$ws.Range("A13").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Rows.Item(13).RowHeight = 60
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# ---- Row 14 becomes just the "Short syllabus:" label (old long Portuguese body removed).
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# ---- Row 15 becomes "Programa:" with the same date text already used in row 8
# (pasted as literal text, not parsed as a date serial), keeping column styles.
$ws.Range("A15").Value = "Programa:"
$ws.Rows.Item(15).RowHeight = 120
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial($xlPasteValues)
$ws.Range("C15").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = 0
$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("C13").Copy()
$ws.Range("C15").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---- Row 16 becomes just the "Syllabus:" label (old long Portuguese body removed).
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# ---- Row 17 becomes just "Avaliação:" with default row height.
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# ---- Row 18 becomes "Método:" label, reusing the professor credential value in B/C.
$ws.Range("A18").Value = "Método:"
$ws.Rows.Item(18).RowHeight = 60
$ws.Range("B18").Value = "5840793 - Sérgio Schneider"
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("C18").Value = "5840793 - Sérgio Schneider"
$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---- Row 19 becomes "Critério:" label with the old "Método:" grading-methods text.
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios."
$ws.Range("C19").Value = "Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios."

# ---- Row 20 becomes "Norma de recuperação:" label with the old "Critério:" weighting text.
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

# ---- Row 21 becomes "Bibliografia:" label with the old recovery-exam text, and its
# row height grows from 60 to 120.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# ---- Row 22 becomes "Requisitos:" only (old bibliography list text removed), default height.
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows.Item(22).AutoFit()

# ---- Row 23 loses its "Requisitos:" label (now on row 22) but keeps the
# requirement value in B/C.
$ws.Range("A23").Clear()
$ws.Range("B23").Value = "LOM3234 -  Óptica Física  (Requisito)`n"
$ws.Range("B19").Copy()
$ws.Range("B23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("C23").Value = "LOM3234 -  Óptica Física  (Requisito)`n"
$ws.Range("C19").Copy()
$ws.Range("C23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Rows.Item(23).RowHeight = 30

# ---- The old row 24 (duplicate requirement value) is removed entirely, shrinking
# the sheet's used range from A1:C24 down to A1:C23.
$ws.Rows.Item(24).Delete()
